$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" values are plain numeric text (e.g. "215.65"). The
# column stores everything as text, so those new values are written with a
# leading quote-prefix ( '215.65 ) -- the standard Excel way to force a
# number-looking entry to stay text -- keeping them text instead of letting
# Excel auto-convert them to numbers.

$ws.Range("D2").Value = "25.863.90"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.638.60"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'215.65"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'0.5052"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.2577"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.06434"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'19.69"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.07790"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.864.82"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.637.27"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'0.5614"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "'63.06"
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "25.889.50"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'194.57"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").Value = "'9.896"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'6.103"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'1.777"
$ws.Range("E25").Value = "  -6.49%  "
$ws.Range("D26").Value = "'140.11"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "'0.1260"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'6.839"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'15.42"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'0.04884"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "'3.295"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "'3.222"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").Value = "'1.569"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "'2.378"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "'0.9036"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'2.580"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'0.5515"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "1.126.01"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'0.01563"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").Value = "'0.9970"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'5.548"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'0.8013"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'98.07"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "1.775.45"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("D47").Value = "'55.42"
$ws.Range("D48").Value = "'0.4261"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").Value = "'7.739"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("D50").Value = "'0.05040"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.45%  "
